$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 6, pushing the "ignore A0..A5" rows down to 12..17
$ws.Range("A6:I11").EntireRow.Insert()

# New data for solo referee mode simulation (rows 6-11)
$ws.Range("A6").Value = "Ref 1 simulated"
$ws.Range("B6").Value = 7
$ws.Range("G6").Value = "ON"
$ws.Range("H6").Value = "refbox/decision"
$ws.Range("I6").Value = "1 good"

$ws.Range("A7").Value = "Ref 2 simulated"
$ws.Range("B7").Value = 7
$ws.Range("G7").Value = "ON"
$ws.Range("H7").Value = "refbox/decision"
$ws.Range("I7").Value = "2 good"

$ws.Range("A8").Value = "Ref 3 simulated"
$ws.Range("B8").Value = 7
$ws.Range("G8").Value = "ON"
$ws.Range("H8").Value = "refbox/decision"
$ws.Range("I8").Value = "3 good"

$ws.Range("A9").Value = "Ref 1 simulated"
$ws.Range("B9").Value = 8
$ws.Range("G9").Value = "ON"
$ws.Range("H9").Value = "refbox/decision"
$ws.Range("I9").Value = "1 bad"

$ws.Range("A10").Value = "Ref 2 simulated"
$ws.Range("B10").Value = 8
$ws.Range("G10").Value = "ON"
$ws.Range("H10").Value = "refbox/decision"
$ws.Range("I10").Value = "2 bad"

$ws.Range("A11").Value = "Ref 3 simulated"
$ws.Range("B11").Value = 8
$ws.Range("G11").Value = "ON"
$ws.Range("H11").Value = "refbox/decision"
$ws.Range("I11").Value = "3 bad"

# Update selection to reflect the new active cell
$ws.Range("I11").Select() | Out-Null
